$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 766
$ws.Range("F4").Value = 90
$ws.Range("F7").Value = 4724
$ws.Range("F8").Value = 350
$ws.Range("F9").Value = 525
$ws.Range("F10").Value = 838
$ws.Range("F11").Value = 791
$ws.Range("F17").Value = 1526
$ws.Range("F18").Value = 1400
$ws.Range("F19").Value = 608
$ws.Range("F21").Value = 163
$ws.Range("F22").Value = 220
$ws.Range("F23").Value = 430
$ws.Range("F24").Value = 89
$ws.Range("F28").Value = 861
$ws.Range("F31").Value = 147
$ws.Range("F37").Value = 496
$ws.Range("F38").Value = 56

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 89

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 227

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 227
$ws.Range("F3").Value = 766
$ws.Range("F5").Value = 90
$ws.Range("F9").Value = 4724
$ws.Range("F10").Value = 350
$ws.Range("F11").Value = 525
$ws.Range("F14").Value = 838
$ws.Range("F15").Value = 791
$ws.Range("F17").Value = 89
$ws.Range("F24").Value = 1526
$ws.Range("F25").Value = 1400
$ws.Range("F26").Value = 608
$ws.Range("F28").Value = 163
$ws.Range("F29").Value = 220
$ws.Range("F31").Value = 430
$ws.Range("F32").Value = 89
$ws.Range("F36").Value = 861
$ws.Range("F39").Value = 147
$ws.Range("F45").Value = 496
$ws.Range("F46").Value = 56
